# Add the new "IP" worksheet (positioned after the existing "logindata" sheet)
# and populate it with the two WBM/MP URLs, matching the author's commit
# ("Added snm pom classes.") which introduced a lookup sheet of IP-based URLs.

$wb = $excel.ActiveWorkbook
$loginSheet = $wb.Worksheets.Item(1)

# New sheet goes right after "logindata" -> becomes sheet index 2 / tab 2 (0-based activeTab=1)
$ipSheet = $wb.Worksheets.Add($null, $loginSheet)
$ipSheet.Name = "IP"

# Populate A2 first, then A1, so the shared-string table gets the same
# insertion order as the authentic edit (wbm before mp).
$ipSheet.Range("A2").Value = "http://10.211.162.111/wbm"
$ipSheet.Range("A1").Value = "http://10.211.162.111/mp"

# Size column A to fit its contents (drops tabSelected off "logindata" and
# puts it on the new, now-active "IP" sheet automatically).
$ipSheet.Columns.Item(1).AutoFit()

# Leave the same selection state captured in the original file.
$ipSheet.Range("B6").Select()
